$wb = $excel.ActiveWorkbook

# --- Fix up Test4's sheet view: it should no longer be the selected tab,
#     and its selection should reflect the sheet's used range (D8:G16),
#     not the stale C5 selection. ---
$ws4 = $wb.Worksheets.Item("Test4")
$ws4.Select() | Out-Null
$ws4.Range("D8:G16").Select() | Out-Null

# --- Add new worksheet "Test5" at the end of the workbook ---
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "Test5"

# --- Populate Test5 with the same table as Test4, shifted down 9 rows and
#     right 2 columns (F17:I25), mirroring the row9/col-E "a" cell lost
#     in the copy. ---
$ws5.Range("F17").Value = "A"
$ws5.Range("G17").Value = "B"
$ws5.Range("H17").Value = "C"
$ws5.Range("I17").Value = "D"

$ws5.Range("H18").Value = "z"

$ws5.Range("G19").Value = "b"
$ws5.Range("H19").Value = "y"
$ws5.Range("I19").Value = 1

$ws5.Range("G20").Value = "c"
$ws5.Range("H20").Value = "x"
$ws5.Range("I20").Value = 2

$ws5.Range("F21").Value = 3
$ws5.Range("G21").Value = "d"
$ws5.Range("H21").Value = "w"
$ws5.Range("I21").Value = 3

$ws5.Range("F22").Value = 4
$ws5.Range("G22").Value = "e"
$ws5.Range("I22").Value = 4

$ws5.Range("F23").Value = 5
$ws5.Range("G23").Value = "f"
$ws5.Range("H23").Value = "v"
$ws5.Range("I23").Value = 5

$ws5.Range("F24").Value = 6
$ws5.Range("G24").Value = "g"
$ws5.Range("H24").Value = "u"

$ws5.Range("G25").Value = "h"

# --- Make Test5 the active sheet/tab, with the selection left just past
#     the used range (I27), matching the target view state. ---
$ws5.Select() | Out-Null
$ws5.Range("I27").Select() | Out-Null
